$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '73.023.81'
$ws.Range('E2').Value = '  +3.79%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.984.10'
$ws.Range('E3').Value = '  +1.95%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '587.83'
$ws.Range('E5').Value = '  +9.48%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '159.05'
$ws.Range('E6').Value = '  +9.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.685'
$ws.Range('E7').Value = '  +0.67%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.36%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.752'
$ws.Range('E9').Value = '  +3.31%  '

$ws.Range('E10').Value = '  +2.04%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '54.00'
$ws.Range('E11').Value = '  +2.28%  '

$ws.Range('E12').Value = '  +2.83%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.87'
$ws.Range('E13').Value = '  +4.41%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.617.65'
$ws.Range('E14').Value = '  +1.47%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.981.22'
$ws.Range('E15').Value = '  +1.60%  '

$ws.Range('E16').Value = '  +10.51%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.10'
$ws.Range('E17').Value = '  +2.52%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.43'
$ws.Range('E18').Value = '  +1.48%  '

$ws.Range('E19').Value = '  +0.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.673.14'
$ws.Range('E20').Value = '  +3.44%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '434.36'
$ws.Range('E21').Value = '  +2.19%  '

$ws.Range('E22').Value = '  +12.54%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '96.41'
$ws.Range('E23').Value = '  +0.72%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.44'
$ws.Range('E24').Value = '  -0.76%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '14.39'
$ws.Range('E25').Value = '  +2.41%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.41'
$ws.Range('E26').Value = '  +23.56%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.12'
$ws.Range('E27').Value = '  +0.43%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.70'
$ws.Range('E28').Value = '  +2.55%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.94'
$ws.Range('E29').Value = '  +1.58%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.54'
$ws.Range('E30').Value = '  +1.46%  '

$ws.Range('E31').Value = '  +6.75%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '13.67'
$ws.Range('E32').Value = '  +3.25%  '

$ws.Range('E33').Value = '  +3.65%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '678.54'
$ws.Range('E34').Value = '  +0.26%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '48.52'
$ws.Range('E35').Value = '  +3.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '69.23'
$ws.Range('E36').Value = '  +7.44%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0882'
$ws.Range('E37').Value = '  +8.18%  '

$ws.Range('E38').Value = '  +2.77%  '

$ws.Range('E39').Value = '  +0.20%  '

$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('B41').Value = 'ThetaToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  -1.75%  '

$ws.Range('B42').Value = 'WEMIXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.34'
$ws.Range('E42').Value = '  +2.32%  '

$ws.Range('E43').Value = '  -0.04%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0488'
$ws.Range('E44').Value = '  +2.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.81'
$ws.Range('E45').Value = '  +12.97%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.150'
$ws.Range('E46').Value = '  +2.24%  '

$ws.Range('E47').Value = '  +0.53%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.42'
$ws.Range('E48').Value = '  +2.69%  '

$ws.Range('E49').Value = '  +3.14%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.42'
$ws.Range('E50').Value = '  +6.73%  '

$ws.Range('E51').Value = '  +8.90%  '
